$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (new quarterly period data),
# shifting the existing D:K data right to F:M.
$ws.Range("D:E").Insert()

# Copy number formats / styles from the (now shifted) old columns F:G into
# the newly inserted D:E columns so the new cells match the existing layout.
$ws.Range("F:G").Copy($ws.Range("D:E"))

# Row 7 / 38 / 80 - Period Ending header dates
$ws.Range("D7").Value = 43373
$ws.Range("E7").Value = 43281
$ws.Range("D38").Value = 43373
$ws.Range("E38").Value = 43281
$ws.Range("D80").Value = 43373
$ws.Range("E80").Value = 43281

# Row 8 - Total Revenue
$ws.Range("D8").Value = 25100
$ws.Range("E8").Value = 8300

# Row 17 - Selling General and Administrative
$ws.Range("D17").Value = 3700
$ws.Range("E17").Value = 1300

# Row 18 - Total Operating Expenses
$ws.Range("D18").Value = 21400
$ws.Range("E18").Value = 7000

# Row 20 - Operating Income or Loss
$ws.Range("D20").Value = -17200
$ws.Range("E20").Value = -5900

# Row 23 - Income Before Tax
$ws.Range("D23").Value = 4200
$ws.Range("E23").Value = 1200

# Row 24 - Income Tax Expense
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 100

# Row 26 - Income After Tax
$ws.Range("D26").Value = 4200
$ws.Range("E26").Value = 1100

# Row 27 - Net Income From Continuing Ops
$ws.Range("D27").Value = 4200
$ws.Range("E27").Value = 1100

# Row 32 - Total Other Income/Expenses Net
$ws.Range("D32").Value = 17200
$ws.Range("E32").Value = 5900

# Row 33 - Net Income
$ws.Range("D33").Value = 4200
$ws.Range("E33").Value = 1100

# Row 35 - Net Income Applicable To Common Shares
$ws.Range("D35").Value = 4200
$ws.Range("E35").Value = 1100

# Row 41 - Cash And Cash Equivalents
$ws.Range("D41").Value = 682500
$ws.Range("E41").Value = 640900

# Row 54 - Total Assets
$ws.Range("D54").Value = 969100
$ws.Range("E54").Value = 944400

# Row 61 - Long Term Debt
$ws.Range("D61").Value = 20600
$ws.Range("E61").Value = 20600

# Row 81 - Net Income (cash flow statement)
$ws.Range("D81").Value = 4200
$ws.Range("E81").Value = 1100
